$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 299
$ws.Range("I2").Value = 98.833336
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 98.833336
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = 14.166664
$ws.Range("N2").Value = -1726
$ws.Range("H5").Value = 91.166664
$ws.Range("I5").Value = 91.166664
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 91.166664
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = 23.833336
$ws.Range("H12").Value = 1178.75
$ws.Range("I12").Value = 572
$ws.Range("J12").Value = 2999
$ws.Range("K12").Value = 572
$ws.Range("L12").Value = 2999
$ws.Range("M12").Value = -402
$ws.Range("N12").Value = -3339
$ws.Range("H62").Value = 6353.88
$ws.Range("I62").Value = 6253
$ws.Range("J62").Value = 6533.222
$ws.Range("K62").Value = 6253
$ws.Range("L62").Value = 6533.222
$ws.Range("M62").Value = -5629
$ws.Range("N62").Value = -7781.222
$ws.Range("H65").Value = 6353.88
$ws.Range("I65").Value = 6253
$ws.Range("J65").Value = 6533.222
$ws.Range("K65").Value = 31265
$ws.Range("L65").Value = 32666.11
$ws.Range("M65").Value = -28145
$ws.Range("N65").Value = -38906.11
$ws.Range("H92").Value = 23810140
$ws.Range("I92").Value = 23810140
$ws.Range("K92").Value = 23810140
$ws.Range("M92").Value = -23808892
$ws.Range("H137").Value = 1711409.5
$ws.Range("I137").Value = 77091.45
$ws.Range("J137").Value = 2528568.5
$ws.Range("K137").Value = 231274.35
$ws.Range("L137").Value = 7585705.5
$ws.Range("M137").Value = -228724.35
$ws.Range("N137").Value = -7590805.5
$ws.Range("H138").Value = 4418.98
$ws.Range("J138").Value = 4499.9893
$ws.Range("L138").Value = 13499.9679
$ws.Range("N138").Value = -23779.9679

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3143.2
$ws.Range("I74").Value = 3010.5715
$ws.Range("K74").Value = 3010.5715
$ws.Range("M74").Value = -2136.5715
$ws.Range("H77").Value = 3143.2
$ws.Range("I77").Value = 3010.5715
$ws.Range("K77").Value = 15052.8575
$ws.Range("M77").Value = -10684.8575
$ws.Range("H97").Value = 780.2083
$ws.Range("I97").Value = 661.95654
$ws.Range("K97").Value = 661.95654
$ws.Range("M97").Value = -165.95654

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").ClearContents()
$ws.Range("N62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").ClearContents()
$ws.Range("N65").Value = 0

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 688.5769
$ws.Range("I16").Value = 645.4375
$ws.Range("J16").Value = 757.6
$ws.Range("K16").Value = 645.4375
$ws.Range("L16").Value = 757.6
$ws.Range("M16").Value = -358.4375
$ws.Range("N16").Value = -1331.6
$ws.Range("H31").Value = 4780.436
$ws.Range("I31").Value = 2998.8
$ws.Range("J31").Value = 5042.4414
$ws.Range("K31").Value = 2998.8
$ws.Range("L31").Value = 5042.4414
$ws.Range("M31").Value = -2703.8
$ws.Range("N31").Value = -5632.4414
$ws.Range("H34").Value = 4780.436
$ws.Range("I34").Value = 2998.8
$ws.Range("J34").Value = 5042.4414
$ws.Range("K34").Value = 2998.8
$ws.Range("L34").Value = 5042.4414
$ws.Range("M34").Value = -2796.8
$ws.Range("N34").Value = -5446.4414
$ws.Range("H113").Value = 688.5769
$ws.Range("I113").Value = 645.4375
$ws.Range("J113").Value = 757.6
$ws.Range("K113").Value = 645.4375
$ws.Range("L113").Value = 757.6
$ws.Range("M113").Value = 1524.5625
$ws.Range("N113").Value = -5097.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 499.75
$ws.Range("I41").Value = 499.75
$ws.Range("K41").Value = 1499.25
$ws.Range("M41").Value = -1161.25
$ws.Range("H68").Value = 2446.4187
$ws.Range("J68").Value = 2633.1428
$ws.Range("L68").Value = 7899.428400000001
$ws.Range("N68").Value = -9521.428400000001
$ws.Range("H71").Value = 2446.4187
$ws.Range("J71").Value = 2633.1428
$ws.Range("L71").Value = 23698.2852
$ws.Range("N71").Value = -31810.2852
$ws.Range("H107").Value = 685.8
$ws.Range("I107").Value = 627.8570999999999
$ws.Range("J107").Value = 1497
$ws.Range("K107").Value = 1883.5713
$ws.Range("L107").Value = 4491
$ws.Range("M107").Value = 36.42870000000016
$ws.Range("N107").Value = -8331
$ws.Range("H134").Value = 6158
$ws.Range("I134").Value = 6158
$ws.Range("K134").Value = 18474
$ws.Range("M134").Value = -13404

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1899
$ws.Range("I107").Value = 2499.5
$ws.Range("J107").Value = 1298.5
$ws.Range("K107").Value = 2499.5
$ws.Range("L107").Value = 1298.5
$ws.Range("M107").Value = -579.5
$ws.Range("N107").Value = -5138.5
$ws.Range("H113").Value = 27232.56
$ws.Range("I113").Value = 4324.3887
$ws.Range("J113").Value = 86139.28999999999
$ws.Range("K113").Value = 4324.3887
$ws.Range("L113").Value = 86139.28999999999
$ws.Range("M113").Value = -2154.3887
$ws.Range("N113").Value = -90479.28999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2407.4285
$ws.Range("I61").Value = 1885.5714
$ws.Range("J61").Value = 2929.2856
$ws.Range("K61").Value = 1885.5714
$ws.Range("L61").Value = 2929.2856
$ws.Range("M61").Value = -1683.5714
$ws.Range("N61").Value = -3333.2856
$ws.Range("H113").Value = 2407.4285
$ws.Range("I113").Value = 1885.5714
$ws.Range("J113").Value = 2929.2856
$ws.Range("K113").Value = 1885.5714
$ws.Range("L113").Value = 2929.2856
$ws.Range("M113").Value = 284.4286
$ws.Range("N113").Value = -7269.2856
$ws.Range("H132").Value = 5695.273
$ws.Range("I132").Value = 5695.273
$ws.Range("K132").Value = 17085.819
$ws.Range("M132").Value = -14555.819
$ws.Range("H136").Value = 6074.276
$ws.Range("I136").Value = 4480.5293
$ws.Range("J136").Value = 8332.083000000001
$ws.Range("K136").Value = 13441.5879
$ws.Range("L136").Value = 24996.249
$ws.Range("M136").Value = -10891.5879
$ws.Range("N136").Value = -30096.249

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 411.3846
$ws.Range("I100").Value = 209.4
$ws.Range("J100").Value = 1084.6666
$ws.Range("K100").Value = 418.8
$ws.Range("L100").Value = 2169.3332
$ws.Range("M100").Value = 122.2
$ws.Range("N100").Value = -3251.3332
$ws.Range("H113").Value = 276.75
$ws.Range("I113").Value = 244.85715
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 734.5714499999999
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1435.42855
$ws.Range("N113").Value = -5840
$ws.Range("H122").Value = 3993.225
$ws.Range("I122").Value = 3992.7812
$ws.Range("K122").Value = 11978.3436
$ws.Range("M122").Value = -9528.3436
